# Replaced backordered part with equivalent.
#
# Row 2 of the Car_Circuitry_BOM sheet (BAT_IN) referenced a backordered
# WAGO terminal block (2601-3102 / footprint 26013102). Swap it for the
# equivalent connector part 2604-3102 (footprint 26043102).
#
# Columns: A=Comment, B=Description, C=Designator, D=Footprint, E=LibRef, F=Quantity
# Leading apostrophes force the values to be stored as text (matching the
# original cells, which were text-formatted/quote-prefixed numeric-looking
# strings) instead of being auto-converted to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Car_Circuitry_BOM")

$ws.Range("A2").Value = "'2604-3102"   # Comment
$ws.Range("B2").Value = "'Connector"   # Description
$ws.Range("D2").Value = "'26043102"    # Footprint
$ws.Range("E2").Value = "'2604-3102"   # LibRef
